$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Swap the worker rows: row 16 (document # + name) and row 17 (document # + name)
# so that TIRSO SAUL ATENCIO ATENCIO (73236943) now appears first (row 16)
# and LILIANA REYES MUÑOZ (33253257) appears second (row 17).
$ws.Range("C16").Value = "73236943"
$ws.Range("D16").Value = "TIRSO SAUL ATENCIO ATENCIO"
$ws.Range("C17").Value = "33253257"
$ws.Range("D17").Value = "LILIANA REYES MUÑOZ"
